# Update the "three-digit number divided by one-digit number" worksheet
# with a newly generated set of division problems.

$d = $word.ActiveDocument

$pairs = @(
    @("321÷8=", "228÷4="),
    @("494÷4=", "165÷4="),
    @("357÷6=", "611÷3="),
    @("380÷8=", "245÷4="),
    @("788÷3=", "154÷5="),
    @("737÷4=", "531÷6="),
    @("155÷2=", "301÷8="),
    @("425÷5=", "891÷2="),
    @("586÷2=", "163÷7="),
    @("139÷2=", "893÷6="),
    @("186÷4=", "436÷2="),
    @("565÷7=", "498÷8="),
    @("979÷3=", "415÷7="),
    @("949÷4=", "216÷4="),
    @("607÷4=", "119÷5="),
    @("491÷5=", "794÷4="),
    @("680÷3=", "627÷6="),
    @("557÷8=", "327÷4="),
    @("409÷7=", "337÷7="),
    @("396÷4=", "780÷7="),
    @("729÷7=", "716÷2="),
    @("620÷8=", "630÷7="),
    @("501÷3=", "562÷9="),
    @("142÷7=", "301÷7="),
    @("624÷6=", "425÷4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
